$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 3; $row -le 23; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H = PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I = LAST UPDATE
    # Leading apostrophe forces Excel to store this as literal text
    # instead of auto-converting the date-like string into a date serial.
    $iCell.Value2 = "'04-Nov-2025"
}
